$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "51.276.45"
$ws.Range("E2").Value = "  +1.26%  "

# Row 3
$ws.Range("D3").Value = "2.963.29"
$ws.Range("E3").Value = "  +2.21%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.26"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.44%  "

# Row 11
$ws.Range("E11").Value = "  +0.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0841"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.85%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.430.94"
$ws.Range("E13").Value = "  +2.71%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.35%  "

# Row 16
$ws.Range("D16").Value = "2.970.09"
$ws.Range("E16").Value = "  +2.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.965"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.43%  "

# Row 18
$ws.Range("D18").Value = "51.295.33"
$ws.Range("E18").Value = "  +0.75%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.33"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.32%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.88"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.75%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0960"
$ws.Range("E22").Value = "  +1.97%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.99%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "260.99"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.69%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.32%  "

# Row 26
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +22.63%  "

# Row 27
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.84%  "

# Row 28
$ws.Range("E28").Value = "  +0.38%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.113"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.23%  "

# Row 30
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.88"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.24%  "

# Row 33
$ws.Range("E33").Value = "  +1.59%  "

# Row 34
$ws.Range("E34").Value = "  -2.09%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.15"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.14%  "

# Row 36
$ws.Range("E36").Value = "  +7.42%  "

# Row 37
$ws.Range("E37").Value = "  +0.27%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.06"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.23"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.84"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.03%  "

# Row 42
$ws.Range("E42").Value = "  +2.85%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.47%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.78%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.290"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +23.99%  "

# Row 46
$ws.Range("E46").Value = "  -0.38%  "

# Row 47
$ws.Range("E47").Value = "  +3.90%  "

# Row 48
$ws.Range("D48").Value = "2.041.56"
$ws.Range("E48").Value = "  +1.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.24"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0344"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +11.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.75%  "

